$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns ("ownTeam", "oppTeam") before the existing column D
# ("batsman"), shifting the old D:I ("batsman".."sr") over to F:K.
$ws.Range("D1:E1").EntireColumn.Insert()

# New header labels for the inserted columns.
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# Fill in the new columns for the existing data row (row 2).
$ws.Range("D2").Value = "Kolkata Knight Riders"
$ws.Range("E2").Value = "Royal Challengers Bangalore"

# Add a new data row (row 3) for the Sharjah match. Force the row to be
# stored as text (matching the rest of the sheet, which keeps numeric-
# looking values such as run counts / strike rate as text) before writing
# the values so figures like "8" and "66.66" aren't coerced into numbers.
$ws.Range("A3:K3").NumberFormat = "@"
$ws.Range("A3").Value = " Sharjah"
$ws.Range("B3").Value = " October 12 2020"
$ws.Range("C3").Value = "RCB won by 82 runs"
$ws.Range("D3").Value = "Kolkata Knight Riders"
$ws.Range("E3").Value = "Royal Challengers Bangalore"
$ws.Range("F3").Value = "Tom Banton "
$ws.Range("G3").Value = "8"
$ws.Range("H3").Value = "12"
$ws.Range("I3").Value = "0"
$ws.Range("J3").Value = "0"
$ws.Range("K3").Value = "66.66"
